# Auto-generated Excel COM-interop script to update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.858.44'
$ws.Range("E2").Value = '  -1.51%  '
$ws.Range("D3").Value = '2.354.64'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'503.92"
$ws.Range("E5").Value = '  -0.26%  '
$ws.Range("D6").Value = "'129.88"
$ws.Range("E6").Value = '  -2.28%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D9").Value = '2.367.98'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  +2.04%  '
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D14").Value = '2.773.90'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").Value = '55.835.91'
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").Value = '2.349.46'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("E19").Value = '  -2.94%  '
$ws.Range("D20").Value = "'310.57"
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = "'4.01"
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("D22").Value = "'6.20"
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("D24").Value = "'65.30"
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("E28").Value = '  -3.30%  '
$ws.Range("D29").Value = "'170.79"
$ws.Range("E29").Value = '  -2.64%  '
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").Value = "'5.72"
$ws.Range("E34").Value = '  -2.71%  '
$ws.Range("E35").Value = '  -5.03%  '
$ws.Range("E36").Value = '  -0.86%  '
$ws.Range("E37").Value = '  -2.29%  '
$ws.Range("D38").Value = "'0.835"
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("E39").Value = '  -4.28%  '
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("E41").Value = '  -3.24%  '
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("D44").Value = "'125.65"
$ws.Range("E44").Value = '  -5.23%  '
$ws.Range("D45").Value = "'0.558"
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").Value = "'240.04"
$ws.Range("E47").Value = '  -2.85%  '
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("D49").Value = "'16.81"
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").Value = "'16.56"
$ws.Range("E51").Value = '  -3.69%  '
